# Adds the key/value pair for z1009 (hyperbola asymptote slope) to Sheet1,
# right after the existing z1008 row, and re-homes the z2xxx / z3xxx blocks
# further down the sheet (matching the gaps used throughout the rest of the
# table) so the layout reads: z1 block (281-289) / gap / z2 block (301-304)
# / gap / z3 block (311-320).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a single new row at 289 for the new z1009 pair. ------------
# This pushes the old z2001:z2004 block (previously 291:294) down to
# 292:295, and the old z3001:z3010 block (previously 301:310) down to
# 302:311.
$ws.Range("A289").EntireRow.Insert()

# --- 2. Insert 9 more blank rows right before the (now shifted) z2 block --
# so it lands on 301:304 (matching the rest of the sheet's "xN01..xN0k"
# blocks always starting on a row ending in 1). This also pushes the z3
# block from 302:311 down to 311:320.
$ws.Range("A292:A300").EntireRow.Insert()

# --- 3. Write the new key/value pair into row 289. -------------------------
$ws.Range("A289").Value = "z1009"
$ws.Range("B289").Value = "쌍곡선의 방정식에서 점근선의 기울기를 구합니다."

# --- 4. Update the view: scroll position + active selection. --------------
$win = $excel.ActiveWindow
$win.ScrollRow = 286
$win.ScrollColumn = 1
$ws.Range("B299").Select()
